$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A5 value (recomputed timestamp with slightly different precision)
$ws.Range("A5").Value = 45865.1668905787

# Add new row 6 data
$ws.Range("A6").NumberFormat = $ws.Range("A5").NumberFormat
$ws.Range("A6").Value = 45865.2086103067
$ws.Range("B6").Value = 2025
$ws.Range("C6").Value = 30
$ws.Range("D6").Value = 13.56
$ws.Range("E6").Value = 89.56999999999999
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 6.61
$ws.Range("H6").Value = "SE"
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = "05:00:23"
